$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric "Importe" values (column H): fix decimal/thousands separators ---
# Force text format first so Excel does not reinterpret these as numbers
# (which would drop the trailing zeros / thousand grouping).
$importeRange = $ws.Range("H2:H149")
$importeRange.NumberFormat = "@"

$ws.Range("H2").Value = "880.00"
$ws.Range("H3").Value = "900.00"
$ws.Range("H4").Value = "21800.00"
$ws.Range("H5").Value = "272.96"
$ws.Range("H6").Value = "354.00"
$ws.Range("H7").Value = "137.00"
$ws.Range("H8").Value = "53914.29"
$ws.Range("H9").Value = "5019.30"
$ws.Range("H10").Value = "482.01"
$ws.Range("H11").Value = "6671.15"
$ws.Range("H12").Value = "7489.76"
$ws.Range("H13").Value = "206.40"
$ws.Range("H14").Value = "259.80"
$ws.Range("H15").Value = "7358.98"
$ws.Range("H16").Value = "427.00"
$ws.Range("H17").Value = "720.00"
$ws.Range("H18").Value = "631.50"
$ws.Range("H19").Value = "43.04"
$ws.Range("H20").Value = "3468.50"
$ws.Range("H21").Value = "10310.00"
$ws.Range("H22").Value = "390.00"
$ws.Range("H23").Value = "11443.00"
$ws.Range("H24").Value = "6538.76"
$ws.Range("H25").Value = "2393.00"
$ws.Range("H26").Value = "101.70"
$ws.Range("H27").Value = "38154.40"
$ws.Range("H28").Value = "23560.00"
$ws.Range("H29").Value = "4289.68"
$ws.Range("H30").Value = "9288.00"
$ws.Range("H31").Value = "2435.00"
$ws.Range("H32").Value = "1138.00"
$ws.Range("H33").Value = "18340.00"
$ws.Range("H34").Value = "17305.90"
$ws.Range("H35").Value = "376.06"
$ws.Range("H36").Value = "49.50"
$ws.Range("H37").Value = "549.50"
$ws.Range("H38").Value = "2178.90"
$ws.Range("H39").Value = "18.00"
$ws.Range("H40").Value = "725.86"
$ws.Range("H41").Value = "1033.00"
$ws.Range("H42").Value = "5826.72"
$ws.Range("H43").Value = "95.00"
$ws.Range("H44").Value = "2313.22"
$ws.Range("H45").Value = "2368.00"
$ws.Range("H46").Value = "7680.00"
$ws.Range("H47").Value = "222.83"
$ws.Range("H48").Value = "98.00"
$ws.Range("H49").Value = "22208.00"
$ws.Range("H50").Value = "3723.00"
$ws.Range("H51").Value = "673.60"
$ws.Range("H52").Value = "270.00"
$ws.Range("H53").Value = "6816.00"
$ws.Range("H54").Value = "1180.00"
$ws.Range("H55").Value = "2192.75"
$ws.Range("H56").Value = "1705.00"
$ws.Range("H57").Value = "290.00"
$ws.Range("H58").Value = "4277.75"
$ws.Range("H59").Value = "2000.00"
$ws.Range("H60").Value = "445.00"
$ws.Range("H61").Value = "913.40"
$ws.Range("H62").Value = "11850.80"
$ws.Range("H63").Value = "2.59"
$ws.Range("H64").Value = "1867.80"
$ws.Range("H65").Value = "3828.90"
$ws.Range("H66").Value = "146.00"
$ws.Range("H67").Value = "1969.00"
$ws.Range("H68").Value = "3092.60"
$ws.Range("H69").Value = "37.00"
$ws.Range("H70").Value = "282.40"
$ws.Range("H71").Value = "373.00"
$ws.Range("H72").Value = "1055.40"
$ws.Range("H73").Value = "105.00"
$ws.Range("H74").Value = "2.18"
$ws.Range("H75").Value = "1848.75"
$ws.Range("H76").Value = "0.20"
$ws.Range("H77").Value = "100170.00"
$ws.Range("H78").Value = "483.50"
$ws.Range("H79").Value = "5.96"
$ws.Range("H80").Value = "0.03"
$ws.Range("H81").Value = "2.40"
$ws.Range("H82").Value = "500.00"
$ws.Range("H83").Value = "1848.80"
$ws.Range("H84").Value = "39.20"
$ws.Range("H85").Value = "71.50"
$ws.Range("H86").Value = "7.65"
$ws.Range("H87").Value = "1702.00"
$ws.Range("H88").Value = "129.31"
$ws.Range("H89").Value = "77.70"
$ws.Range("H90").Value = "1260.22"
$ws.Range("H91").Value = "1870.00"
$ws.Range("H92").Value = "240.00"
$ws.Range("H93").Value = "580.00"
$ws.Range("H94").Value = "1858.00"
$ws.Range("H95").Value = "564.00"
$ws.Range("H96").Value = "5972.00"
$ws.Range("H97").Value = "1250.00"
$ws.Range("H98").Value = "400.00"
$ws.Range("H99").Value = "1928.02"
$ws.Range("H100").Value = "85.00"
$ws.Range("H101").Value = "331.50"
$ws.Range("H102").Value = "500.00"
$ws.Range("H103").Value = "1000.00"
$ws.Range("H104").Value = "2100.00"
$ws.Range("H105").Value = "5776.55"
$ws.Range("H106").Value = "290.00"
$ws.Range("H107").Value = "250.00"
$ws.Range("H108").Value = "450.00"
$ws.Range("H109").Value = "6110.20"
$ws.Range("H110").Value = "200.00"
$ws.Range("H111").Value = "1800.00"
$ws.Range("H112").Value = "350.00"
$ws.Range("H113").Value = "2000.00"
$ws.Range("H114").Value = "1700.00"
$ws.Range("H115").Value = "120.00"
$ws.Range("H116").Value = "8734.97"
$ws.Range("H117").Value = "456.00"
$ws.Range("H118").Value = "167.00"
$ws.Range("H119").Value = "310.00"
$ws.Range("H120").Value = "380.00"
$ws.Range("H121").Value = "5250.00"
$ws.Range("H122").Value = "220.00"
$ws.Range("H123").Value = "602.57"
$ws.Range("H124").Value = "221.00"
$ws.Range("H125").Value = "156.00"
$ws.Range("H126").Value = "2280.00"
$ws.Range("H127").Value = "400.00"
$ws.Range("H128").Value = "60.00"
$ws.Range("H129").Value = "46.60"
$ws.Range("H130").Value = "400.00"
$ws.Range("H131").Value = "68.60"
$ws.Range("H132").Value = "46.21"
$ws.Range("H133").Value = "58345.45"
$ws.Range("H134").Value = "210.00"
$ws.Range("H135").Value = "1983.60"
$ws.Range("H136").Value = "17652.90"
$ws.Range("H137").Value = "1672.50"
$ws.Range("H138").Value = "80.96"
$ws.Range("H139").Value = "0.61"
$ws.Range("H140").Value = "266.70"
$ws.Range("H141").Value = "2557.12"
$ws.Range("H142").Value = "1624.43"
$ws.Range("H143").Value = "6800.00"
$ws.Range("H144").Value = "1265762.32"
$ws.Range("H145").Value = "289264.79"
$ws.Range("H146").Value = "44000.00"
$ws.Range("H147").Value = "64600.00"
$ws.Range("H148").Value = "1000.00"
$ws.Range("H149").Value = "860.00"

# Restore the original (default) cell style now that the text is locked in
$importeRange.Style = "Normal"

# --- Razon social / Nombre Fantasia (columns E/F): normalize separators ---
$ws.Range("E60").Value = "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"
$ws.Range("F60").Value = "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"
$ws.Range("E68").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E70").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F70").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E71").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E118").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F118").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E119").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
